# Automated map update (2025-09-30 06:56:14)
# - "General" sheet: fix row 3 (OT pending load / traspaso propio) + append 4 new rows (381-384)
# - "PEBCOM" sheet: fix row 3 (same case, mirrored on this sheet)
# - "Sin_Asignar" sheet: append the same 4 new rows (62-65)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Row 3 correction shared by "General" and "PEBCOM" (same underlying case)
# ---------------------------------------------------------------------------
$row3Sheets = @("General", "PEBCOM")
foreach ($sheetName in $row3Sheets) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Cells.Item(3, 5).NumberFormat = "@"      # E3 - OT
    $ws.Cells.Item(3, 5).Value = "Pendiente de Carga"

    $ws.Cells.Item(3, 7).NumberFormat = "@"      # G3 - Estado
    $ws.Cells.Item(3, 7).Value = "Pendiente"

    $ws.Cells.Item(3, 8).NumberFormat = "@"      # H3 - Observaciones
    $ws.Cells.Item(3, 8).Value = "Desmontar columna ya traspasaron nodo"
}

# ---------------------------------------------------------------------------
# 2) New rows appended to "General" (rows 381-384) and "Sin_Asignar" (62-65)
# ---------------------------------------------------------------------------
$newRows = @(
    @{ A="-620"; B="9/29/2025"; C="Luis Viale 3098";      D="11"; E="810056875"; F="Sin Asignar"; G="Pendiente"; H="picada"; I=1; J="Cambio"; K="Sin equipos"; L="Terminal"; M=-58.477413; N=-34.620772; O="Devoto";     P="Capital Norte"; Q="NRA-M"; R="Fuera de Poligono OVL" },
    @{ A="-621"; B="9/29/2025"; C="Tres Arroyos 2911";    D="11"; E="810056868"; F="Sin Asignar"; G="Pendiente"; H="Picada"; I=1; J="Cambio"; K="Sin equipos"; L="Pasante";  M=-58.476877; N=-34.617525; O="Devoto";     P="Capital Norte"; Q="NRA-M"; R="Fuera de Poligono OVL" },
    @{ A="-622"; B="9/29/2025"; C="Mariano Acha 2271";    D="12"; E="810056867"; F="Sin Asignar"; G="Pendiente"; H="Picada"; I=1; J="Cambio"; K="Sin equipos"; L="Pasante";  M=-58.477338; N=-34.571921; O="Colegiales"; P="Capital Norte"; Q="ATH-J"; R="Fuera de Poligono OVL" },
    @{ A="-623"; B="9/29/2025"; C="Mosconi 3368";         D="11"; E="810061513"; F="Sin Asignar"; G="Pendiente"; H="Picada"; I=1; J="Cambio"; K="Sin equipos"; L="Pasante";  M=-58.508377; N=-34.590137; O="Paternal";   P="Capital Norte"; Q="PUE-N"; R="ARATO-25058.PO.2PUE" }
)

# Text columns (left as strings, forced with NumberFormat "@" so numeric-looking
# values like "11", "810056875" or "-620" are not silently coerced to numbers)
$textCols = @{
    1 = "A"; 2 = "B"; 3 = "C"; 4 = "D"; 5 = "E"; 6 = "F"; 7 = "G"; 8 = "H";
    10 = "J"; 11 = "K"; 12 = "L"; 15 = "O"; 16 = "P"; 17 = "Q"; 18 = "R"
}
# Numeric columns: 9 = I (Attachments), 13 = M (Coordenada_X), 14 = N (Coordenada_Y)

function Append-Rows($sheetName, $startRow) {
    $ws = $wb.Worksheets.Item($sheetName)
    $r = $startRow
    foreach ($row in $newRows) {
        foreach ($col in $textCols.Keys) {
            $key = $textCols[$col]
            $cell = $ws.Cells.Item($r, $col)
            $cell.NumberFormat = "@"
            $cell.Value = $row[$key]
        }
        $ws.Cells.Item($r, 9).Value = $row["I"]
        $ws.Cells.Item($r, 13).Value = $row["M"]
        $ws.Cells.Item($r, 14).Value = $row["N"]
        $r++
    }
}

Append-Rows "General" 381
Append-Rows "Sin_Asignar" 62
